# "update code after merge"
#
# The sheet used to have a leading "STT" (row-number) column:
#   A: STT   B: Email              C: Password
#   A2: 1    B2: hieu002@gmail.com C2: blue299   (hyperlink lived on B2)
#
# After the merge the STT column is gone and everything shifted one
# column to the left, with the mailto: hyperlink now following its data
# to column A:
#   A: Email              B: Password
#   A2: hieu002@gmail.com B2: blue299            (hyperlink now on A2)
#
# Also bring the sheet's base/hyperlink font in line with the rest of the
# merged workbook (Calibri -> Arial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-theme the Normal / Hyperlink cell styles to Arial.
$wb.Styles.Item("Normal").Font.Name = "Arial"
$wb.Styles.Item("Hyperlink").Font.Name = "Arial"

# The hyperlink collection doesn't auto-follow a column delete, so drop it
# before shifting the data and re-create it afterwards on the new cell.
$ws.Hyperlinks.Delete()

# Remove the obsolete "STT" column; B/C/D slide left into A/B/C.
$ws.Columns.Item(1).Delete()

# Re-attach the mailto: link to its data, which now lives in A2, and make
# sure it keeps the Hyperlink look.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:hieu002@gmail.com")
$ws.Range("A2").Style = "Hyperlink"
